$d = $word.ActiveDocument
$vbVT = [char]11

# -------------------------------------------------------------------------
# 1. Replace the three <w:br/> runs in the title paragraph with " | " runs,
#    while keeping each text segment in its own distinct run (matching how
#    Word itself would preserve run boundaries instead of recomputing /
#    merging the whole paragraph's runs, which is what a naive in-place
#    Find/Replace across the paragraph would otherwise do).
#
#    Technique: temporarily split the title paragraph at each break so the
#    break character is isolated (alone) in its own paragraph, replace the
#    isolated break text (safe - nothing to merge with), then re-join all
#    the paragraphs by deleting only the paragraph marks in between (this
#    does not trigger run recombination the way editing text in place does).
# -------------------------------------------------------------------------

# Locate the three break characters inside paragraph 1 and split the
# paragraph right after each one.
$titlePar = $d.Paragraphs(1)
$breakPositions = @()
$rng = $titlePar.Range.Duplicate()
$rng.Find.ClearFormatting()
while ($rng.Find.Execute($vbVT, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)) {
    $breakPositions += $rng.Start
    $rng.Collapse(0)  # collapse to the end, continue searching forward
}

# Split after each break, processing from the last break to the first so
# earlier offsets remain valid.
for ($i = $breakPositions.Count - 1; $i -ge 0; $i--) {
    $pos = $breakPositions[$i] + 1
    $d.Range($pos, $pos).InsertParagraphAfter() | Out-Null
}

# Now split right before each break too, isolating it into its own
# single-character paragraph. Re-find each break from the last paragraph
# to the first (paragraph indices shift as we insert, so always re-derive
# positions via Find instead of relying on stale offsets).
for ($i = 3; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $pr = $p.Range.Duplicate()
    $pr.Find.ClearFormatting()
    $pr.Find.Execute($vbVT, $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    $d.Range($pr.Start, $pr.Start).InsertParagraphBefore() | Out-Null
}

# Paragraphs 2, 4 and 6 now each contain exactly one break character and
# nothing else - replace each with " | ".
foreach ($i in 2, 4, 6) {
    $p = $d.Paragraphs($i)
    $pr = $p.Range.Duplicate()
    $pr.Find.ClearFormatting()
    $pr.Find.Replacement.ClearFormatting()
    $pr.Find.Execute($vbVT, $false, $false, $false, $false, $false, $true, 1, $false, " | ", 2) | Out-Null
}

# Re-join the 7 temporary paragraphs back into a single paragraph by
# deleting just the paragraph marks between them (not the run text).
for ($k = 0; $k -lt 6; $k++) {
    $p1 = $d.Paragraphs(1)
    $markPos = $p1.Range.End - 1
    $d.Range($markPos, $markPos + 1).Delete() | Out-Null
}

# -------------------------------------------------------------------------
# 2. Delete the empty paragraph, the "Author: ..." paragraph, the
#    "Location: ..." paragraph, and the following empty paragraph.
# -------------------------------------------------------------------------
$rStart = $d.Paragraphs(2).Range.Start
$rEnd = $d.Paragraphs(5).Range.End
$d.Range($rStart, $rEnd).Delete() | Out-Null

# -------------------------------------------------------------------------
# 3. Delete the trailing empty paragraph right before the sectPr. Deleting
#    the very last paragraph's own range leaves its paragraph mark behind
#    (Word always keeps a final mark), so also consume the mark of the
#    paragraph before it.
# -------------------------------------------------------------------------
$n = $d.Paragraphs.Count
$lastPar = $d.Paragraphs($n)
if ($lastPar.Range.Text -eq [char]13) {
    $prevEnd = $d.Paragraphs($n - 1).Range.End
    $lastEnd = $lastPar.Range.End
    $d.Range($prevEnd - 1, $lastEnd).Delete() | Out-Null
}

$d.Save()
